# Updated symbol list on Sun Jan 15 11:46:35 UTC 2023 with GitHub Actions
# Applies updated price/volume values to the cryptos worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    # Prefix with an apostrophe so Excel stores the numeric-looking text
    # (e.g. "296.91", "-2.02%") verbatim as a string instead of coercing
    # it to a floating point number / percentage.
    $c.Value = "'" + $value
    # Re-apply the default "Normal" style so the transient quote-prefix
    # formatting that Excel applies to text-forced cells doesn't stick
    # around on the cell (matches original, unstyled data cells).
    $c.Style = "Normal"
}

Set-TextValue "D2" "296.91"
Set-TextValue "E2" "-2.02%"
Set-TextValue "D3" "31.26"
Set-TextValue "E3" "-1.56%"
Set-TextValue "D4" "5.109"
Set-TextValue "E4" "-1.91%"
Set-TextValue "D5" "0.07336"
Set-TextValue "E5" "-0.39%"
Set-TextValue "D6" "7.722"
Set-TextValue "E6" "-1.49%"
Set-TextValue "D7" "1.667"
Set-TextValue "E7" "10.62%"
Set-TextValue "D8" "3.729"
Set-TextValue "E8" "-0.18%"
Set-TextValue "D9" "0.9200"
Set-TextValue "E9" "1.31%"
Set-TextValue "D10" "0.1679"
Set-TextValue "E10" "-0.15%"
Set-TextValue "D11" "0.07034"
Set-TextValue "E11" "-6.43%"
Set-TextValue "D12" "0.08091"
Set-TextValue "E12" "1.86%"
Set-TextValue "D13" "0.02992"
Set-TextValue "E13" "0.75%"
Set-TextValue "D14" "0.09904"
Set-TextValue "D15" "0.001495"
Set-TextValue "E15" "0.51%"
Set-TextValue "D16" "0.006166"
Set-TextValue "E16" "0.13%"
Set-TextValue "D17" "3.446"
Set-TextValue "E17" "-0.71%"
Set-TextValue "E18" "-0.08%"
Set-TextValue "E19" "-1.98%"
Set-TextValue "D20" "0.1332"
Set-TextValue "E20" "0.64%"
Set-TextValue "E21" "1.16%"
Set-TextValue "D22" "0.04638"
Set-TextValue "E22" "2.55%"
Set-TextValue "D23" "0.1549"
Set-TextValue "E23" "-4.32%"
Set-TextValue "D24" "0.001214"
Set-TextValue "E24" "-0.26%"
Set-TextValue "D26" "0.0001297"
Set-TextValue "E26" "0.00%"
Set-TextValue "E27" "7.72%"
Set-TextValue "D39" "0.01701"
Set-TextValue "E39" "1.78%"
Set-TextValue "D40" "0.04422"
Set-TextValue "E40" "-1.40%"
Set-TextValue "D41" "0.007171"
Set-TextValue "E41" "-0.34%"
Set-TextValue "D42" "0.1329"
Set-TextValue "E42" "-1.07%"
Set-TextValue "D43" "0.002136"
Set-TextValue "E44" "-13.29%"
Set-TextValue "D45" "0.00005983"
Set-TextValue "E45" "-1.26%"
Set-TextValue "D47" "1.894"
Set-TextValue "E47" "0.28%"
